$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price and Volume columns remain text so numeric-looking values
# (e.g. "303.12", "2.290") keep their exact formatting/trailing zeros
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "23.426.93"
$ws.Range("E2").Value = "  +0.66%  "

$ws.Range("D3").Value = "1.639.78"
$ws.Range("E3").Value = "  +1.03%  "

$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").Value = "1.002"
$ws.Range("E5").Value = "  -0.08%  "

$ws.Range("D6").Value = "303.12"

$ws.Range("D7").Value = "0.3820"
$ws.Range("E7").Value = "  +1.52%  "

$ws.Range("D8").Value = "52.03"
$ws.Range("E8").Value = "  -0.18%  "

$ws.Range("D9").Value = "0.3602"
$ws.Range("E9").Value = "  -0.38%  "

$ws.Range("D10").Value = "0.08211"
$ws.Range("E10").Value = "  +1.68%  "

$ws.Range("D11").Value = "1.229"
$ws.Range("E11").Value = "  -0.08%  "

$ws.Range("D12").Value = "1.003"
$ws.Range("E12").Value = "  -0.14%  "

$ws.Range("D13").Value = "22.49"
$ws.Range("E13").Value = "  -0.48%  "

$ws.Range("D14").Value = "6.441"
$ws.Range("E14").Value = "  -1.56%  "

$ws.Range("D15").Value = "7.334"
$ws.Range("E15").Value = "  +1.76%  "

$ws.Range("D16").Value = "0.00001236"
$ws.Range("E16").Value = "  -0.44%  "

$ws.Range("D17").Value = "1.637.59"
$ws.Range("E17").Value = "  +0.82%  "

$ws.Range("D18").Value = "95.17"
$ws.Range("E18").Value = "  +1.93%  "

$ws.Range("E19").Value = "  +0.73%  "

$ws.Range("D20").Value = "6.590"
$ws.Range("E20").Value = "  +2.17%  "

$ws.Range("D21").Value = "17.46"
$ws.Range("E21").Value = "  -2.46%  "

$ws.Range("D22").Value = "1.002"
$ws.Range("E22").Value = "  -0.18%  "

$ws.Range("D23").Value = "12.50"
$ws.Range("E23").Value = "  -1.50%  "

$ws.Range("D24").Value = "23.441.54"
$ws.Range("E24").Value = "  +0.64%  "

$ws.Range("D25").Value = "2.538"
$ws.Range("E25").Value = "  +4.75%  "

$ws.Range("D26").Value = "3.062"
$ws.Range("E26").Value = "  -4.29%  "

$ws.Range("D27").Value = "21.15"
$ws.Range("E27").Value = "  +0.32%  "

$ws.Range("D28").Value = "151.82"
$ws.Range("E28").Value = "  +1.65%  "

$ws.Range("D29").Value = "5.277"
$ws.Range("E29").Value = "  -0.17%  "

$ws.Range("D30").Value = "134.25"
$ws.Range("E30").Value = "  -0.31%  "

$ws.Range("D31").Value = "1.819.63"
$ws.Range("E31").Value = "  +0.70%  "

$ws.Range("D32").Value = "1.085"
$ws.Range("E32").Value = "  +14.59%  "

$ws.Range("D33").Value = "2.153"
$ws.Range("E33").Value = "  -6.46%  "

$ws.Range("D34").Value = "6.496"
$ws.Range("E34").Value = "  -3.85%  "

$ws.Range("D35").Value = "11.48"
$ws.Range("E35").Value = "  +5.13%  "

$ws.Range("D36").Value = "0.02766"
$ws.Range("E36").Value = "  -2.08%  "

$ws.Range("D37").Value = "0.2507"
$ws.Range("E37").Value = "  -0.67%  "

$ws.Range("D38").Value = "0.08771"
$ws.Range("E38").Value = "  -0.44%  "

$ws.Range("D39").Value = "0.07024"
$ws.Range("E39").Value = "  -1.04%  "

$ws.Range("D40").Value = "5.955"
$ws.Range("E40").Value = "  -2.39%  "

$ws.Range("D41").Value = "0.7019"
$ws.Range("E41").Value = "  -0.05%  "

$ws.Range("D42").Value = "1.346"
$ws.Range("E42").Value = "  -0.99%  "

$ws.Range("D43").Value = "12.23"
$ws.Range("E43").Value = "  -0.54%  "

$ws.Range("E44").Value = "  -3.26%  "

$ws.Range("D45").Value = "0.6501"
$ws.Range("E45").Value = "  +0.94%  "

$ws.Range("D46").Value = "1.001"
$ws.Range("E46").Value = "  -0.08%  "

$ws.Range("D47").Value = "2.290"
$ws.Range("E47").Value = "  -0.87%  "

$ws.Range("D48").Value = "3.959"
$ws.Range("E48").Value = "  -0.52%  "

$ws.Range("D49").Value = "0.07969"
$ws.Range("E49").Value = "  +0.07%  "

$ws.Range("D50").Value = "128.82"
$ws.Range("E50").Value = "  +2.24%  "

$ws.Range("D51").Value = "1.192"
$ws.Range("E51").Value = "  -0.98%  "
